$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns (G:K) for the "meta" statistics, shifting the
# existing "arrecadado_*" ... "maior_ano" columns from G:V to L:AA.
$ws.Range("G1:K1").EntireColumn.Insert()

# The inserted header cells (row 1) already inherit the bold/centered
# "header" style (s=4) from the row, so only the text needs to be set.
$ws.Cells.Item(1, 7).Value2  = "meta"
$ws.Cells.Item(1, 8).Value2  = "meta_avg"
$ws.Cells.Item(1, 9).Value2  = "meta_std"
$ws.Cells.Item(1, 10).Value2 = "meta_min"
$ws.Cells.Item(1, 11).Value2 = "meta_max"

# The inserted data cells (rows 2:6) inherit the percentage style from
# column F; copy the currency style (R$ #,##0.00, s=3) that is used by
# the other monetary columns so the new "meta" columns match.
$ws.Range("L2").Copy()
$ws.Range("G2:K6").PasteSpecial(-4122)

# Fill in the new "meta" statistics values.
$ws.Cells.Item(2, 7).Value2  = 721610.3061912227
$ws.Cells.Item(2, 8).Value2  = 10458.12037958294
$ws.Cells.Item(2, 9).Value2  = 11144.2267578863
$ws.Cells.Item(2, 10).Value2 = 44.33046360042423
$ws.Cells.Item(2, 11).Value2 = 50590.198657868

$ws.Cells.Item(3, 7).Value2  = 5883940.636230236
$ws.Cells.Item(3, 8).Value2  = 13372.59235506872
$ws.Cells.Item(3, 9).Value2  = 19267.96260047285
$ws.Cells.Item(3, 10).Value2 = 23.98859826184044
$ws.Cells.Item(3, 11).Value2 = 147790.8327903106

$ws.Cells.Item(4, 7).Value2  = 1712986.472842461
$ws.Cells.Item(4, 8).Value2  = 9732.877686604894
$ws.Cells.Item(4, 9).Value2  = 10102.88946115519
$ws.Cells.Item(4, 10).Value2 = 46.55761904502517
$ws.Cells.Item(4, 11).Value2 = 83151.82469725677

$ws.Cells.Item(5, 7).Value2  = 7150010.825257363
$ws.Cells.Item(5, 8).Value2  = 10347.33838676898
$ws.Cells.Item(5, 9).Value2  = 16064.05218382809
$ws.Cells.Item(5, 10).Value2 = 12.04441558726698
$ws.Cells.Item(5, 11).Value2 = 198811.9434626772

$ws.Cells.Item(6, 7).Value2  = 131168.4623975197
$ws.Cells.Item(6, 8).Value2  = 18738.35177107424
$ws.Cells.Item(6, 9).Value2  = 19781.31029827062
$ws.Cells.Item(6, 10).Value2 = 2420.445520432476
$ws.Cells.Item(6, 11).Value2 = 54319.48382898097
